# Auto-generated edit script replicating the Anima_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7627.3213
$ws.Range("I19").Value = 399.92307
$ws.Range("J19").Value = 13891.066
$ws.Range("K19").Value = 399.92307
$ws.Range("L19").Value = 13891.066
$ws.Range("M19").Value = -224.92307
$ws.Range("N19").Value = -14241.066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1369.1666
$ws.Range("J127").Value = 1543
$ws.Range("L127").Value = 4629
$ws.Range("N127").Value = -14549

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3670.386
$ws.Range("I138").Value = 4144.3
$ws.Range("J138").Value = 3569.5532
$ws.Range("K138").Value = 12432.9
$ws.Range("L138").Value = 10708.6596
$ws.Range("M138").Value = -7292.900000000001
$ws.Range("N138").Value = -20988.6596

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 559.2727
$ws.Range("I2").Value = 573.0714
$ws.Range("J2").Value = 535.125
$ws.Range("K2").Value = 573.0714
$ws.Range("L2").Value = 535.125
$ws.Range("M2").Value = -460.0714
$ws.Range("N2").Value = -761.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10871986
$ws.Range("I74").Value = 1467.1538
$ws.Range("J74").Value = 25003662
$ws.Range("K74").Value = 1467.1538
$ws.Range("L74").Value = 25003662
$ws.Range("M74").Value = -593.1538
$ws.Range("N74").Value = -25005410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10871986
$ws.Range("I77").Value = 1467.1538
$ws.Range("J77").Value = 25003662
$ws.Range("K77").Value = 7335.769
$ws.Range("L77").Value = 125018310
$ws.Range("M77").Value = -2967.769
$ws.Range("N77").Value = -125027046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 559.2727
$ws.Range("I116").Value = 573.0714
$ws.Range("J116").Value = 535.125
$ws.Range("K116").Value = 573.0714
$ws.Range("L116").Value = 535.125
$ws.Range("M116").Value = 1720.9286
$ws.Range("N116").Value = -5123.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 559.2727
$ws.Range("I3").Value = 573.0714
$ws.Range("J3").Value = 535.125
$ws.Range("K3").Value = 573.0714
$ws.Range("L3").Value = 535.125
$ws.Range("M3").Value = -459.0714
$ws.Range("N3").Value = -763.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3999.5
$ws.Range("I105").Value = 3999.5
$ws.Range("K105").Value = 3999.5
$ws.Range("M105").Value = -2252.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 67390.8
$ws.Range("I107").Value = 111413.664
$ws.Range("J107").Value = 1356.5
$ws.Range("K107").Value = 111413.664
$ws.Range("L107").Value = 1356.5
$ws.Range("M107").Value = -109493.664
$ws.Range("N107").Value = -5196.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4889.0605
$ws.Range("I31").Value = 1258.7391
$ws.Range("J31").Value = 6830.8604
$ws.Range("K31").Value = 1258.7391
$ws.Range("L31").Value = 6830.8604
$ws.Range("M31").Value = -963.7391
$ws.Range("N31").Value = -7420.8604

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4889.0605
$ws.Range("I34").Value = 1258.7391
$ws.Range("J34").Value = 6830.8604
$ws.Range("K34").Value = 1258.7391
$ws.Range("L34").Value = 6830.8604
$ws.Range("M34").Value = -1056.7391
$ws.Range("N34").Value = -7234.8604

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1995.8334
$ws.Range("I99").Value = 1567.6666
$ws.Range("J99").Value = 2167.1
$ws.Range("K99").Value = 1567.6666
$ws.Range("L99").Value = 2167.1
$ws.Range("M99").Value = -69.66660000000002
$ws.Range("N99").Value = -5163.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 886.3125
$ws.Range("I105").Value = 866
$ws.Range("J105").Value = 920.1667
$ws.Range("K105").Value = 866
$ws.Range("L105").Value = 920.1667
$ws.Range("M105").Value = 881
$ws.Range("N105").Value = -4414.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1995.8334
$ws.Range("I126").Value = 1567.6666
$ws.Range("J126").Value = 2167.1
$ws.Range("K126").Value = 4702.9998
$ws.Range("L126").Value = 6501.299999999999
$ws.Range("M126").Value = -2232.9998
$ws.Range("N126").Value = -11441.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2328.04
$ws.Range("I132").Value = 2136.8572
$ws.Range("J132").Value = 2571.3635
$ws.Range("K132").Value = 6410.571599999999
$ws.Range("L132").Value = 7714.0905
$ws.Range("M132").Value = -3880.571599999999
$ws.Range("N132").Value = -12774.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 709.8246
$ws.Range("J5").Value = 1015.3333
$ws.Range("L5").Value = 3045.9999
$ws.Range("N5").Value = -3269.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 800
$ws.Range("J92").Value = 800
$ws.Range("L92").Value = 2400
$ws.Range("N92").Value = -4896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 709.8246
$ws.Range("J135").Value = 1015.3333
$ws.Range("L135").Value = 9137.9997
$ws.Range("N135").Value = -14207.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 30261
$ws.Range("J39").Value = 30261
$ws.Range("L39").Value = 30261
$ws.Range("N39").Value = -31325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 19166.666
$ws.Range("I122").Value = 36000
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 108000
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -105550
$ws.Range("N122").Value = -11899.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2760.25
$ws.Range("I126").Value = 2077
$ws.Range("J126").Value = 3319.2727
$ws.Range("K126").Value = 6231
$ws.Range("L126").Value = 9957.8181
$ws.Range("M126").Value = -3761
$ws.Range("N126").Value = -14897.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5441.8
$ws.Range("I7").Value = 5433.9473
$ws.Range("J7").Value = 5466.6665
$ws.Range("K7").Value = 5433.9473
$ws.Range("L7").Value = 5466.6665
$ws.Range("M7").Value = -5321.9473
$ws.Range("N7").Value = -5690.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3231
$ws.Range("I40").Value = 3038.75
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3038.75
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2902.75
$ws.Range("N40").Value = -4272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2748.2856
$ws.Range("I61").Value = 2225.0908
$ws.Range("J61").Value = 4666.6665
$ws.Range("K61").Value = 2225.0908
$ws.Range("L61").Value = 4666.6665
$ws.Range("M61").Value = -2023.0908
$ws.Range("N61").Value = -5070.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2748.2856
$ws.Range("I113").Value = 2225.0908
$ws.Range("J113").Value = 4666.6665
$ws.Range("K113").Value = 2225.0908
$ws.Range("L113").Value = 4666.6665
$ws.Range("M113").Value = -55.09079999999994
$ws.Range("N113").Value = -9006.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5441.8
$ws.Range("I126").Value = 5433.9473
$ws.Range("J126").Value = 5466.6665
$ws.Range("K126").Value = 16301.8419
$ws.Range("L126").Value = 16399.9995
$ws.Range("M126").Value = -13831.8419
$ws.Range("N126").Value = -21339.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1146.6177
$ws.Range("I122").Value = 1182.9259
$ws.Range("K122").Value = 3548.7777
$ws.Range("M122").Value = -1098.7777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1331.3846
$ws.Range("I126").Value = 1275.6666
$ws.Range("K126").Value = 3826.9998
$ws.Range("M126").Value = -1356.9998

Write-Output "Applied Anima_Profits updates"
